# Update TestData.xlsx with latest ParaBank test data
#
# - ParaBank_RegistrationForm (sheet5): the "account created" message now
#   wraps across three lines, and the row's Status flips from FAIL to PASS.
# - ParaBank_LoginForm (sheet6) becomes inactive; ParaBank_RegistrationForm
#   becomes the active / selected sheet.

$wb = $excel.ActiveWorkbook

$wsReg = $wb.Worksheets.Item("ParaBank_RegistrationForm")

# Multi-line "account created" success message (embedded line breaks) with
# word-wrap turned on for the cell.
$msg = "Your account was created `nsuccessfully. You are now `nlogged in."
$cellMsg = $wsReg.Range("M2")
$cellMsg.Value = $msg
$cellMsg.WrapText = $true

# Status for that row moves from FAIL to PASS.
$wsReg.Range("N2").Value = "PASS"

# Make the registration form the active sheet/tab (was the login form).
$wsReg.Activate()
